$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New "pull correct / push correct" helper table in columns K:O.
# Cell writes are ordered to reproduce the shared-string insertion
# order of the authored workbook (new unique strings first-seen in
# this exact sequence).
# ------------------------------------------------------------------

# Row 1 - left sensor header (bold, like the other section headers)
$ws.Range("M1").Value = "Left Sensor Maximum Dist from Left Wall"
$ws.Range("M1").Font.Bold = $true

# Row 4 - right sensor header (bold)
$ws.Range("M4").Value = "Right Sensor Maximum Dist from Right Wall"
$ws.Range("M4").Font.Bold = $true

# Row 2 - range either side of the left-sensor max distance
$ws.Range("M2").Value = 360
$ws.Range("N2").Value = ">>"
$ws.Range("O2").Value = 440

# Row 7 - center labels
$ws.Range("M7").Value = "Left Center"
$ws.Range("N7").Value = "Right Center"

# Row 10 - off-center labels + edge-of-frame note
$ws.Range("M10").Value = "L Off Ctr"
$ws.Range("N10").Value = "R Off Ctr"
$ws.Range("L11").Value = "4.75 cm from left"
$ws.Range("L12").Value = "4.75 cm from right"
$ws.Range("L10").Value = "From Edge of Frame to Wall"
$ws.Range("K10").Value = "AC"
$ws.Range("L13").Value = "x cm from left"
$ws.Range("L14").Value = "x cm from right"

# Row 5 - range either side of the right-sensor max distance
$ws.Range("M5").Value = 280
$ws.Range("N5").Value = ">>"
$ws.Range("O5").Value = 380

# Row 8 - push/pull correct values
$ws.Range("M8").Value = 550
$ws.Range("N8").Value = 410

# Rows 11-14 - AC group numbering + measurements
$ws.Range("K11").Value = 1
$ws.Range("M11").Value = 580
$ws.Range("N11").Value = 400

$ws.Range("K12").Value = 1
$ws.Range("M12").Value = 520
$ws.Range("N12").Value = 450

$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 440

$ws.Range("K14").Value = 2
$ws.Range("N14").Value = 380

# Highlight fills that call out the push-correct / pull-correct pair
# (theme Accent1, lighter 40% / lighter 60% - closest reproducible
# approximation of the authored tinted theme fills). The darker fill
# is applied first so it claims the lower fill/style index, matching
# the authored style table order.
$ws.Range("M13").Interior.Color = 14136213
$ws.Range("N14").Interior.Color = 14136213
$ws.Range("M11").Interior.Color = 14994616
$ws.Range("N12").Interior.Color = 14994616

# ------------------------------------------------------------------
# Column widths for the new helper columns.
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 5.3072916666666666
$ws.Columns.Item(12).ColumnWidth = 17.451822916666668
$ws.Columns.Item(13).ColumnWidth = 11.451822916666666
$ws.Columns.Item(14).ColumnWidth = 11.592447916666666

# ------------------------------------------------------------------
# Reposition the chart so it no longer overlaps the new K:O data.
# ------------------------------------------------------------------
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = $chartObj.Left + 267
$chartObj.Top = $chartObj.Top + 24

# ------------------------------------------------------------------
# View state - selection moved, matching the authored workbook.
# ------------------------------------------------------------------
$ws.Range("N19").Select()
